$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently sits right
#    after the title (Heading1) paragraph.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new paragraph just before the final paragraph (the one
#    that holds the old "Prompt: ..." image-generation text) containing
#    the bolded title line, matching the exact run layout used
#    elsewhere in the document (an empty leading run + a formatted run).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$count2 = $d.Paragraphs.Count
$titlePara = $d.Paragraphs.Item($count2 - 1)
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Combinator Free - Innovative Cascading Mechanism</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara.Range.InsertXML($titleXml)

# 3. Swap the old image-prompt text in the (now) final paragraph for the
#    meta-description sentence, keeping its existing (italic) run
#    formatting untouched.
$oldText = 'Prompt: Create a feature image for "Fruit Combinator" in a cartoon style featuring a happy Maya warrior with glasses Description: The feature image for "Fruit Combinator" will feature a cartoon-style design, with a happy Maya warrior wearing glasses. The image will be colorful and playful, with bold lines and bright colors. The Maya warrior will be holding a bowl of fruit, with different fruit combinations spilling out around them, showcasing the game''s theme of progressive cascades and ways to win. In the background of the image, we will see a jungle scene, with vines and plants winding around the edges of the frame. The overall feel of the image will be upbeat and fun, capturing the excitement and thrill of playing the game. The Maya warrior, with their quirky glasses and wide grin, will convey a sense of joy and excitement, encouraging players to dive into the world of "Fruit Combinator."'
$newText = 'Read our review of Fruit Combinator - an online slot game with an innovative cascade mechanism. Play for free before betting with real money.'
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
